$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: best_params strings
$ws.Range("B2").Value = "{'alpha': 0.1, 'max_iter': 1000}"
$ws.Range("C2").Value = "{'alpha': 0.001, 'max_iter': 1000}"
$ws.Range("D2").Value = "{'alpha': 0.001, 'l1_ratio': 0.25, 'max_iter': 1000}"
$ws.Range("E2").Value = "{'C': 1, 'gamma': 1}"
$ws.Range("K2").Value = "{'activation': 'relu', 'b_random_vec_range': [0, 10], 'lam': 1, 'n_layer': 16, 'n_nodes': 256, 'random_seed': 358, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3: rmse
$ws.Range("B3").Value = 0.07365191382727573
$ws.Range("C3").Value = 0.09482543562865793
$ws.Range("D3").Value = 0.0822952872556722
$ws.Range("E3").Value = 0.07614074535858557
$ws.Range("F3").Value = 0.05290137422467665
$ws.Range("G3").Value = 0.05019867482937072
$ws.Range("H3").Value = 0.1015212474423683
$ws.Range("I3").Value = 0.04719998031908251
$ws.Range("J3").Value = 0.05070881903069058
$ws.Range("K3").Value = 0.02673389891972606

# Row 4: r2
$ws.Range("B4").Value = 0.9049519609820754
$ws.Range("C4").Value = 0.8423433234817181
$ws.Range("D4").Value = 0.881720028899452
$ws.Range("E4").Value = 0.8988318188252457
$ws.Range("F4").Value = 0.9483827196779414
$ws.Range("G4").Value = 0.953396796924791
$ws.Range("H4").Value = 0.8230625889627975
$ws.Range("I4").Value = 0.9595943643042307
$ws.Range("J4").Value = 0.9542264983691717
$ws.Range("K4").Value = 0.986653353629201

# Row 5: mape
$ws.Range("B5").Value = 13.20094430097138
$ws.Range("C5").Value = 17.97304476803705
$ws.Range("D5").Value = 15.27312564475146
$ws.Range("E5").Value = 16.13094892244744
$ws.Range("F5").Value = 6.427331223082751
$ws.Range("G5").Value = 6.371517571339483
$ws.Range("H5").Value = 19.40265142933292
$ws.Range("I5").Value = 6.296702272859105
$ws.Range("J5").Value = 7.175632936858507
$ws.Range("K5").Value = 3.774573562720013
